# fix para exceles con rows vacias
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Explicacion")

# Explicacion (Sheet2): convert date-serial cells in column I to text strings
# (processed first so the new shared string "04/04/2022" gets the lower index)
$ws2.Range("I16").Value = "04/04/2022"
$ws2.Range("I17").Value = "04/04/2022"

# Sheet1: convert date-serial cells in column I to text strings
$ws1.Range("I2").Value = "22/05/2022"
$ws1.Range("I3").Value = "22/05/2022"
$ws1.Range("I3").Font.Underline = $true

# Add an extra (empty) row below on Explicacion, touching I18 so an
# empty-but-styled row exists (created after I3's style so the cellXfs
# index ordering matches)
$ws2.Range("I18").Font.Underline = $true
$ws2.Range("I18").Value = ""
